$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old layout (A1 and A3:B6 block) before rebuilding it.
$ws.Range("A1:B6").ClearContents()

# Row 1: "week" label and a numeric week value.
$ws.Range("B1").Value = "week"
$ws.Range("C1").Value = 42

# Row 2: title moved down from row 1.
$ws.Range("A2").Value = "Microcontroller Team 1"

# Row 4-7: team members, each name+surname merged into a single cell.
$ws.Range("A4").Value = "Muhammad Amjad" + " " + "Bin Abdul Malik"
$ws.Range("A5").Value = "Muhammad Iqbal" + " " + "Bin Mohd Fauzi"
$ws.Range("A6").Value = "Muhammad Farid Izwan" + " " + "Bin Mohamad Shabri"
$ws.Range("C6").Value = "Coding for Arduino Circuit Pedestrian &Car traffic light"
$ws.Range("A7").Value = "Muhammad Amirul Hakimi " + " " + "Bin Zaprunnizam"

# Row 8: new team member.
$ws.Range("A8").Value = "Patrick Stephen "

$ws.Range("C2").Select()
